$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2917716402565462
$ws.Range("C2").Value = 10.34677158129881
$ws.Range("D2").Value = 22.3905356188092
$ws.Range("E2").Value = 2195978.878461985
$ws.Range("G2").Value = 2196011.907540825
